$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Managed three member group ..." -> "Managed three members group ..."
#    The canonical XML shows the inserted "s" landing in its own run
#    (classic "typed one character mid-sentence" run split), so we
#    locate the insertion point and insert the character rather than
#    doing a blind Find&Replace of the whole sentence.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("three member", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter("s")
# Force the newly-inserted run to carry its own explicit direct
# formatting (splitting it from its neighbour) the same way Word
# does when a character is typed in the middle of a run.
$r.Font.Underline = 1
$r.Font.Underline = 0

# ------------------------------------------------------------------
# 2) Drop the stray "HR ACR(Oct'2020-Nov'2020)" that had been glued
#    onto the Technologies line.
# ------------------------------------------------------------------
$d.Content.Find.Execute("PHP (Laravel), Docker, JQuery, MySQL.HR ACR(Oct" + [char]8217 + "2020-Nov" + [char]8217 + "2020)", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "PHP (Laravel), Docker, JQuery, MySQL.", 2)

# ------------------------------------------------------------------
# 3) Tiny inline "cursor" shape resized slightly.
# ------------------------------------------------------------------
$shp = $d.InlineShapes.Item(1)
$shp.Width = 7620 / 12700
$shp.Height = 26035 / 12700

# ------------------------------------------------------------------
# 4) Skills table column-width rounding tweaks (table 1, columns 7 & 8).
# ------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$t1.Columns.Item(7).Width = 994 / 20
$t1.Columns.Item(8).Width = 900 / 20

# ------------------------------------------------------------------
# 5) Education table column-width rounding tweaks (table 2, columns 2 & 3).
# ------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$t2.Columns.Item(2).Width = 3149 / 20
$t2.Columns.Item(3).Width = 3336 / 20
